$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the "http://www.toronto.ca/TWONTO#" prefix from the superclass
# identifiers in column A (rows 3-7), splitting the long URI strings into
# new, shorter shared-string entries while leaving column B untouched.
$ws.Range("A3").Value = "instrumentation"
$ws.Range("A4").Value = "air_duct_segment"
$ws.Range("A5").Value = "cable_segment"
$ws.Range("A6").Value = "electrical_panel_or_cabinet"
$ws.Range("A7").Value = "instrument_gauge_or_display"

# Move the active selection off the table, matching the author's saved view.
$ws.Range("A9").Select()
